$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.451.49"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "1.799.45"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'224.34"
$ws.Range("E5").Value = "  -1.88%  "
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'39.03"
$ws.Range("E8").Value = "  +6.47%  "
$ws.Range("E9").Value = "  -4.71%  "
$ws.Range("E10").Value = "  -4.74%  "
$ws.Range("D11").Value = "'0.0982"
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").Value = "2.058.56"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "'10.88"
$ws.Range("E13").Value = "  -5.33%  "
$ws.Range("D14").Value = "1.798.70"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "34.409.47"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "'0.627"
$ws.Range("E16").Value = "  -4.38%  "
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("E18").Value = "  -3.07%  "
$ws.Range("D19").Value = "'238.79"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("D20").Value = "0.0₃0763"
$ws.Range("E20").Value = "  -4.18%  "
$ws.Range("D21").Value = "'11.06"
$ws.Range("E21").Value = "  -4.98%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "'4.07"
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("E24").Value = "  -3.78%  "
$ws.Range("D25").Value = "'170.36"
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("D27").Value = "'7.65"
$ws.Range("E27").Value = "  -4.72%  "
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("E33").Value = "  -4.84%  "
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").Value = "'0.638"
$ws.Range("E35").Value = "  -5.04%  "
$ws.Range("D37").Value = "1.301.86"
$ws.Range("E37").Value = "  -7.12%  "
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("D39").Value = "'2.30"
$ws.Range("E39").Value = "  -6.54%  "
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("D42").Value = "'81.52"
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("D43").Value = "'2.80"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("D44").Value = "'0.942"
$ws.Range("D45").Value = "'13.98"
$ws.Range("E45").Value = "  +3.10%  "
$ws.Range("E46").Value = "  +4.52%  "
$ws.Range("D47").Value = "1.958.93"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("E48").Value = "  -6.00%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").Value = "'101.69"
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("D51").Value = "'0.0612"
$ws.Range("E51").Value = "  -0.62%  "
